$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.222.39'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.858.34'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7102'
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.01'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9991'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08007'
$ws.Range('E8').Value = '  +4.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3033'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08190'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '1.855.33'
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.176'
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7048'
$ws.Range('E14').Value = '  -3.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.76'
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('D16').Value = '29.182.96'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.835'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007877'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.28'
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '238.29'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9979'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.083.18'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9996'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.442'
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.41'
$ws.Range('E25').Value = '  +0.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.932'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1445'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.10'
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.933'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.431'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.480'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.371'
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05213'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7103'
$ws.Range('E36').Value = '  +1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9969'
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.671'
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01860'
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.724'
$ws.Range('E40').Value = '  +1.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9290'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('D42').Value = '1.129.06'
$ws.Range('E42').Value = '  +4.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4263'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '70.37'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.851'
$ws.Range('E45').Value = '  -2.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9989'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.92'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5339'
$ws.Range('E48').Value = '  -4.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.766'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.171'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('D51').Value = '1.974.77'
$ws.Range('E51').Value = '  -1.01%  '
